# Auto-generated edit script applying the cryptos.xlsx diff
# Updates Coin/Link/Price/Volume(1h) columns for rows 2-51 per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.823.26'
$ws.Range("E2").Value = '  -3.25%  '
$ws.Range("D3").Value = '2.650.40'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.597'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.30%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.578'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.53%  '
$ws.Range("D13").Value = '3.052.70'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.107'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").Value = '2.653.10'
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.923'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.61%  '
$ws.Range("D18").Value = '45.783.01'
$ws.Range("E18").Value = '  -4.28%  '
$ws.Range("E19").Value = '  -2.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '279.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '30.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.30%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.37'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.70%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.24%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.70'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.22%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '155.73'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0838'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.08%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.90%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.121'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.55%  '
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.94%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.03'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0326'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.77%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.28%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.95'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.50%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.152.67'
$ws.Range("E45").Value = '  +2.99%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '93.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.43%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '111.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.02%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.912.10'
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.200'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.34%  '
